# Insert two new data rows (709 and 710) into the daily price log sheet,
# pushing all existing rows from 709 downward by two positions (709->711,
# 710->712, ..., 772->774). The sheet's used range grows from A1:T772 to
# A1:T774.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block (formats get inherited
# from the row above, matching the existing column D date style).
$ws.Range("A709:T710").EntireRow.Insert()

# --- New row 709 ---
$ws.Range("A709").Value = 4
$ws.Range("B709").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C709").Value = "Los Lagos"
$ws.Range("D709").Value = 45013
$ws.Range("E709").Value = 10
$ws.Range("F709").Value = "Fruta"
$ws.Range("G709").Value = 100102
$ws.Range("H709").Value = "Cítricos"
$ws.Range("I709").Value = 100102003
$ws.Range("J709").Value = "Limón"
$ws.Range("K709").Value = "Sin especificar"
$ws.Range("L709").Value = "1a plateado"
$ws.Range("M709").Value = 1200
$ws.Range("N709").Value = 30000
$ws.Range("O709").Value = 31000
$ws.Range("P709").Value = 30500
$ws.Range("Q709").Value = "`$/malla 18 kilos"
$ws.Range("R709").Value = "Región de O'Higgins"
$ws.Range("S709").Value = 1694
$ws.Range("T709").Value = 18

# --- New row 710 ---
$ws.Range("A710").Value = 4
$ws.Range("B710").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C710").Value = "Los Lagos"
$ws.Range("D710").Value = 45013
$ws.Range("E710").Value = 10
$ws.Range("F710").Value = "Fruta"
$ws.Range("G710").Value = 100102
$ws.Range("H710").Value = "Cítricos"
$ws.Range("I710").Value = 100102003
$ws.Range("J710").Value = "Limón"
$ws.Range("K710").Value = "Sin especificar"
$ws.Range("L710").Value = "2a plateado"
$ws.Range("M710").Value = 600
$ws.Range("N710").Value = 27000
$ws.Range("O710").Value = 27000
$ws.Range("P710").Value = 27000
$ws.Range("Q710").Value = "`$/malla 18 kilos"
$ws.Range("R710").Value = "Región de O'Higgins"
$ws.Range("S710").Value = 1500
$ws.Range("T710").Value = 18
